$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellPlain($cell, $kind, $val) {
  if ($kind -eq "absent") {
    $cell.ClearContents() | Out-Null
  } elseif ($kind -eq "str") {
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.ClearFormats() | Out-Null
  } else {
    $cell.Value2 = $val
  }
}

# ---- Row 2 ----
$row = 2
$data = @(
  @{Col=1; Kind='num'; Val=258538},
  @{Col=2; Kind='num'; Val=92804},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='NT'},
  @{Col=5; Kind='num'; Val=782},
  @{Col=6; Kind='str'; Val='Skirmossa'},
  @{Col=7; Kind='str'; Val='Hookeria lucens'},
  @{Col=8; Kind='str'; Val='(Hedw.) Sm.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, strax N om Balketorp, Dls'},
  @{Col=17; Kind='num'; Val=333081.3534867804},
  @{Col=18; Kind='num'; Val=6498342.744522936},
  @{Col=19; Kind='num'; Val=10},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='1999-05-15'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='1999-05-15'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='str'; Val='Tämligen rikligt längs ca 150 m av bäcken. Datum osäkert.'},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='str'; Val='I skuggig blandskog'},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='str'; Val='På fuktig jord i bäckkant'},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Leif Appelgren'},
  @{Col=50; Kind='str'; Val='Leif Appelgren'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 3 ----
$row = 3
$data = @(
  @{Col=1; Kind='num'; Val=2199338},
  @{Col=2; Kind='num'; Val=108193},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=219711},
  @{Col=6; Kind='str'; Val='Sårläka'},
  @{Col=7; Kind='str'; Val='Sanicula europaea'},
  @{Col=8; Kind='str'; Val='L.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Balketorp, 400 m NO om, Dls'},
  @{Col=17; Kind='num'; Val=333273.3609825537},
  @{Col=18; Kind='num'; Val=6498504.655406407},
  @{Col=19; Kind='num'; Val=50},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='1979-04-22'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='1979-04-22'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Kjell Eriksson'},
  @{Col=50; Kind='str'; Val='Kjell Eriksson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 4 ----
$row = 4
$data = @(
  @{Col=1; Kind='num'; Val=2192816},
  @{Col=2; Kind='num'; Val=104489},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=219686},
  @{Col=6; Kind='str'; Val='Vätteros'},
  @{Col=7; Kind='str'; Val='Lathraea squamaria'},
  @{Col=8; Kind='str'; Val='L.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Balketorp, 400 m NO om, Dls'},
  @{Col=17; Kind='num'; Val=333273.3609825537},
  @{Col=18; Kind='num'; Val=6498504.655406407},
  @{Col=19; Kind='num'; Val=50},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='1982-05-19'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='1982-05-19'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Kjell Eriksson'},
  @{Col=50; Kind='str'; Val='Kjell Eriksson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 5 ----
$row = 5
$data = @(
  @{Col=1; Kind='num'; Val=2202431},
  @{Col=2; Kind='num'; Val=108194},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=219711},
  @{Col=6; Kind='str'; Val='Sårläka'},
  @{Col=7; Kind='str'; Val='Sanicula europaea'},
  @{Col=8; Kind='str'; Val='L.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Balketorp, 600 m N-NNV om, Dls'},
  @{Col=17; Kind='num'; Val=332860.5715872086},
  @{Col=18; Kind='num'; Val=6498729.504914329},
  @{Col=19; Kind='num'; Val=50},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='1983-06-10'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='1983-06-10'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Kjell Eriksson'},
  @{Col=50; Kind='str'; Val='Kjell Eriksson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 7 ----
$row = 7
$data = @(
  @{Col=1; Kind='num'; Val=69173232},
  @{Col=2; Kind='num'; Val=90655},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='NT'},
  @{Col=5; Kind='num'; Val=788},
  @{Col=6; Kind='str'; Val='Gul taggsvamp'},
  @{Col=7; Kind='str'; Val='Hydnellum geogenium'},
  @{Col=8; Kind='str'; Val='(Fr.) Banker'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, Dls'},
  @{Col=17; Kind='num'; Val=333304.9626084958},
  @{Col=18; Kind='num'; Val=6498681.251649193},
  @{Col=19; Kind='num'; Val=10},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2017-08-27'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2017-08-27'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='str'; Val='Äldre barrdominerad skog utmed bäck med delvis kvillande lopp'},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Henrik Weibull'},
  @{Col=50; Kind='str'; Val='Henrik Weibull'},
  @{Col=51; Kind='str'; Val='Åtgärdsprogram för mossor i Västra Götalands län'}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 8 ----
$row = 8
$data = @(
  @{Col=1; Kind='num'; Val=69173229},
  @{Col=2; Kind='num'; Val=92864},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='VU'},
  @{Col=5; Kind='num'; Val=815},
  @{Col=6; Kind='str'; Val='Stor skogsbäckmossa'},
  @{Col=7; Kind='str'; Val='Hygrohypnum subeugyrium'},
  @{Col=8; Kind='str'; Val='(Renauld & Cardot) Broth.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, Dls'},
  @{Col=17; Kind='num'; Val=333243.743277251},
  @{Col=18; Kind='num'; Val=6498692.219576385},
  @{Col=19; Kind='num'; Val=10},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2017-08-27'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2017-08-27'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='str'; Val='Äldre barrdominerad skog utmed bäck med delvis kvillande lopp'},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Henrik Weibull'},
  @{Col=50; Kind='str'; Val='Henrik Weibull'},
  @{Col=51; Kind='str'; Val='Åtgärdsprogram för mossor i Västra Götalands län'}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 9 ----
$row = 9
$data = @(
  @{Col=1; Kind='num'; Val=69173228},
  @{Col=2; Kind='num'; Val=108194},
  @{Col=3; Kind='str'; Val='Godkänd baserat på observatörens uppgifter'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=219711},
  @{Col=6; Kind='str'; Val='Sårläka'},
  @{Col=7; Kind='str'; Val='Sanicula europaea'},
  @{Col=8; Kind='str'; Val='L.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, Dls'},
  @{Col=17; Kind='num'; Val=333243.743277251},
  @{Col=18; Kind='num'; Val=6498692.219576385},
  @{Col=19; Kind='num'; Val=10},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2017-08-27'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2017-08-27'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='str'; Val='Äldre barrdominerad skog utmed bäck med delvis kvillande lopp'},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Henrik Weibull'},
  @{Col=50; Kind='str'; Val='Henrik Weibull'},
  @{Col=51; Kind='str'; Val='Åtgärdsprogram för mossor i Västra Götalands län'}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 10 ----
$row = 10
$data = @(
  @{Col=1; Kind='num'; Val=16044867},
  @{Col=2; Kind='num'; Val=90671},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='NT'},
  @{Col=5; Kind='num'; Val=4368},
  @{Col=6; Kind='str'; Val='Dofttaggsvamp'},
  @{Col=7; Kind='str'; Val='Hydnellum suaveolens'},
  @{Col=8; Kind='str'; Val='(Scop.:Fr.) P. Karst.'},
  @{Col=9; Kind='str'; Val='10'},
  @{Col=10; Kind='str'; Val='fruktkroppar'},
  @{Col=11; Kind='str'; Val=''},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='str'; Val=''},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Buxåskullen, 675 m NNO Balketorp, Dls'},
  @{Col=17; Kind='num'; Val=333364.6442548583},
  @{Col=18; Kind='num'; Val=6498805.502037385},
  @{Col=19; Kind='num'; Val=25},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2014-07-01'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2014-07-01'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='str'; Val=''},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='str'; Val='Granskog'},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='str'; Val='Kjell Eriksson'},
  @{Col=44; Kind='str'; Val='F1402'},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='str'; Val='Rolf-Göran Carlsson'},
  @{Col=48; Kind='str'; Val='2021'},
  @{Col=49; Kind='str'; Val='Kjell Eriksson'},
  @{Col=50; Kind='str'; Val='Kjell Eriksson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 11 ----
$row = 11
$data = @(
  @{Col=1; Kind='num'; Val=96142344},
  @{Col=2; Kind='num'; Val=90319},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=4769},
  @{Col=6; Kind='str'; Val='Svavelriska'},
  @{Col=7; Kind='str'; Val='Lactarius scrobiculatus'},
  @{Col=8; Kind='str'; Val='(Scop.:Fr.) Fr.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, Dls'},
  @{Col=17; Kind='num'; Val=333265.6097484134},
  @{Col=18; Kind='num'; Val=6498726.768297401},
  @{Col=19; Kind='num'; Val=5},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2021-09-14'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2021-09-14'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Anton Larsson'},
  @{Col=50; Kind='str'; Val='Anton Larsson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 12 ----
$row = 12
$data = @(
  @{Col=1; Kind='num'; Val=96142343},
  @{Col=2; Kind='num'; Val=90319},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=4769},
  @{Col=6; Kind='str'; Val='Svavelriska'},
  @{Col=7; Kind='str'; Val='Lactarius scrobiculatus'},
  @{Col=8; Kind='str'; Val='(Scop.:Fr.) Fr.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelidsbäcken, Dls'},
  @{Col=17; Kind='num'; Val=333251.1676608387},
  @{Col=18; Kind='num'; Val=6498743.562633296},
  @{Col=19; Kind='num'; Val=5},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2021-09-14'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2021-09-14'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Anton Larsson'},
  @{Col=50; Kind='str'; Val='Anton Larsson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

# ---- Row 13 ----
$row = 13
$data = @(
  @{Col=1; Kind='num'; Val=96142350},
  @{Col=2; Kind='num'; Val=90653},
  @{Col=3; Kind='str'; Val='Ovaliderad'},
  @{Col=4; Kind='str'; Val='LC'},
  @{Col=5; Kind='num'; Val=4364},
  @{Col=6; Kind='str'; Val='Dropptaggsvamp'},
  @{Col=7; Kind='str'; Val='Hydnellum ferrugineum'},
  @{Col=8; Kind='str'; Val='(Fr.:Fr.) P. Karst.'},
  @{Col=9; Kind='str'; Val=''},
  @{Col=10; Kind='absent'; Val=$null},
  @{Col=11; Kind='absent'; Val=$null},
  @{Col=12; Kind='absent'; Val=$null},
  @{Col=13; Kind='absent'; Val=$null},
  @{Col=14; Kind='absent'; Val=$null},
  @{Col=15; Kind='absent'; Val=$null},
  @{Col=16; Kind='str'; Val='Muggelid, Dls'},
  @{Col=17; Kind='num'; Val=333288.6459826281},
  @{Col=18; Kind='num'; Val=6498947.551675561},
  @{Col=19; Kind='num'; Val=5},
  @{Col=20; Kind='str'; Val='Västra Götaland'},
  @{Col=21; Kind='str'; Val='Färgelanda'},
  @{Col=22; Kind='str'; Val='Dalsland'},
  @{Col=23; Kind='str'; Val='Färgelanda'},
  @{Col=24; Kind='absent'; Val=$null},
  @{Col=25; Kind='str'; Val='2021-09-14'},
  @{Col=26; Kind='str'; Val='00:00'},
  @{Col=27; Kind='str'; Val='2021-09-14'},
  @{Col=28; Kind='str'; Val='00:00'},
  @{Col=29; Kind='absent'; Val=$null},
  @{Col=30; Kind='bool'; Val=$false},
  @{Col=31; Kind='bool'; Val=$false},
  @{Col=32; Kind='absent'; Val=$null},
  @{Col=33; Kind='bool'; Val=$false},
  @{Col=34; Kind='absent'; Val=$null},
  @{Col=35; Kind='absent'; Val=$null},
  @{Col=36; Kind='absent'; Val=$null},
  @{Col=37; Kind='absent'; Val=$null},
  @{Col=38; Kind='absent'; Val=$null},
  @{Col=39; Kind='absent'; Val=$null},
  @{Col=40; Kind='absent'; Val=$null},
  @{Col=41; Kind='absent'; Val=$null},
  @{Col=42; Kind='absent'; Val=$null},
  @{Col=43; Kind='absent'; Val=$null},
  @{Col=44; Kind='absent'; Val=$null},
  @{Col=45; Kind='absent'; Val=$null},
  @{Col=46; Kind='str'; Val=''},
  @{Col=47; Kind='absent'; Val=$null},
  @{Col=48; Kind='absent'; Val=$null},
  @{Col=49; Kind='str'; Val='Anton Larsson'},
  @{Col=50; Kind='str'; Val='Anton Larsson'},
  @{Col=51; Kind='str'; Val=''}
)
foreach ($item in $data) {
  $cell = $ws.Cells.Item($row, $item.Col)
  Set-CellPlain $cell $item.Kind $item.Val
}

Write-Host "Row permutation applied."
